$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cell values (row-wise, matching target diff)
$ws.Range("B2").Value = -19
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = -1
$ws.Range("E2").Value = -3
$ws.Range("F2").Value = 12
$ws.Range("G2").Value = 7
$ws.Range("I2").Value = -2
$ws.Range("H3").Value = -2
$ws.Range("I3").Value = -5
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = -1
$ws.Range("H5").Value = -1
$ws.Range("I5").Value = 1
$ws.Range("B7").Value = -13
$ws.Range("C7").Value = -4
$ws.Range("D7").Value = -3
$ws.Range("E7").Value = -1
$ws.Range("F7").Value = 8
$ws.Range("G7").Value = -20
$ws.Range("H7").Value = -5
$ws.Range("I7").Value = -2
$ws.Range("F8").Value = 8
$ws.Range("G8").Value = -22
$ws.Range("B9").Value = -20
$ws.Range("C9").Value = -7
$ws.Range("D9").Value = -2
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = -10
$ws.Range("G9").Value = 12
$ws.Range("H9").Value = 8
$ws.Range("I9").Value = -3
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("B11").Value = -36
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = -1
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 7
$ws.Range("G11").Value = -18

# Clear cells that no longer have content in the target sheet
$ws.Range("B3:G3").ClearContents()
$ws.Range("B4:I4").ClearContents()
$ws.Range("B5:C5").ClearContents()
$ws.Range("F5:G5").ClearContents()
$ws.Range("B6:I6").ClearContents()
$ws.Range("B8:E8").ClearContents()
$ws.Range("H8:I8").ClearContents()
$ws.Range("B10:C10").ClearContents()
$ws.Range("F10:I10").ClearContents()
$ws.Range("H11:I11").ClearContents()

